# feat: add MUC1 to genomic targeting
#
# Append a new data row for the MUC1 gene to the genomic_targeting sheet:
#   gene_symbol = MUC1
#   targeting   = TRUE
#   notes       = ADTKD-MUC1 - complete genomic targeting

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row right after the current data (row 14 -> 15)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "MUC1"
$ws.Cells.Item($newRow, 2).Value = $true
$ws.Cells.Item($newRow, 3).Value = "ADTKD-MUC1 - complete genomic targeting"
